# Update the "cryptos" price/volume table with refreshed market data.
# Cells whose new text could be parsed by Excel as a plain number (e.g. "308.97",
# "12.20", "92.20") are written with a leading apostrophe so Excel keeps them as
# literal text (matching the original inlineStr/text storage and preserving
# trailing zeros / decimal formatting) instead of silently converting them to
# numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.233.96'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '2.271.40'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'" + '308.97'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = "'" + '97.13'
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").Value = "'" + '0.527'
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = "'" + '0.489'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").Value = "'" + '35.12'
$ws.Range("E10").Value = '  -3.27%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = "'" + '6.82'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '2.624.90'
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = "'" + '14.68'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '2.272.94'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = "'" + '0.789'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").Value = '42.113.25'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = "'" + '12.20'
$ws.Range("E19").Value = '  -4.57%  '
$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").Value = "'" + '5.97'
$ws.Range("D22").Value = "'" + '67.61'
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("D23").Value = "'" + '236.30'
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").Value = "'" + '2.59'
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").Value = "'" + '0.999'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = "'" + '23.54'
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = "'" + '37.14'
$ws.Range("E28").Value = '  -1.58%  '
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = "'" + '163.73'
$ws.Range("E31").Value = '  +1.54%  '
$ws.Range("D32").Value = "'" + '5.24'
$ws.Range("E32").Value = '  -1.70%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = "'" + '3.09'
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = "'" + '17.59'
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'" + '0.0734'
$ws.Range("E36").Value = '  -2.70%  '
$ws.Range("D37").Value = "'" + '2.37'
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("E39").Value = '  -3.88%  '
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D41").Value = "'" + '4.16'
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("D42").Value = "'" + '2.27'
$ws.Range("E42").Value = '  -6.53%  '
$ws.Range("D43").Value = '1.949.29'
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("D45").Value = "'" + '18.82'
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("D46").Value = "'" + '2.96'
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("D47").Value = "'" + '9.77'
$ws.Range("E47").Value = '  -4.75%  '
$ws.Range("D48").Value = "'" + '53.91'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '2.497.28'
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").Value = "'" + '92.20'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = "'" + '71.50'
$ws.Range("E51").Value = '  -2.14%  '
